$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Convert M16:N21 and M23:N26 into shared formulas (fill-down pattern) ---
$ws.Range("M16:N21").Formula = "=I12/18.5"
$ws.Range("M23:N26").Formula = "=I21/18.5"

# --- 2) New data block: rows 69-80, columns B (2) and C (3) ---
$bcVals = @(
    @(-0.32432432432432434, 0.097297297297297303),
    @(-0.25945945945945947, 0.097297297297297303),
    @(-0.25945945945945947, 0.16216216216216217),
    @(-0.22702702702702704, 0.16216216216216217),
    @(-0.22702702702702704, 0.097297297297297303),
    @(-0.16216216216216217, 0.097297297297297303),
    @(-0.16216216216216217, 0.064864864864864868),
    @(-0.22702702702702704, 0.064864864864864868),
    @(-0.22702702702702704, 0),
    @(-0.25945945945945947, 0),
    @(-0.25945945945945947, 0.064864864864864868),
    @(-0.32432432432432434, 0.064864864864864868)
)
for ($i = 0; $i -lt $bcVals.Length; $i++) {
    $row = 69 + $i
    $pair = $bcVals[$i]
    $ws.Cells.Item($row, 2).Value2 = $pair[0]
    $ws.Cells.Item($row, 3).Value2 = $pair[1]
}

# --- 3) New data block: rows 69-77, columns E (5) and F (6) ---
$efVals1 = @(
    @(-0.32432432432432434, 0.097297297297297303),
    @(-0.25945945945945947, 0.097297297297297303),
    @(-0.22702702702702704, 0.097297297297297303),
    @(-0.16216216216216217, 0.097297297297297303),
    @(-0.16216216216216217, 0.064864864864864868),
    @(-0.22702702702702704, 0.064864864864864868),
    @(-0.25945945945945947, 0.064864864864864868),
    @(-0.32432432432432434, 0.064864864864864868),
    @(-0.32432432432432434, 0.097297297297297303)
)
for ($i = 0; $i -lt $efVals1.Length; $i++) {
    $row = 69 + $i
    $pair = $efVals1[$i]
    $ws.Cells.Item($row, 5).Value2 = $pair[0]
    $ws.Cells.Item($row, 6).Value2 = $pair[1]
}

# --- 4) New data block: rows 79-87, columns E (5) and F (6) ---
$efVals2 = @(
    @(-0.25945945945945947, 0.097297297297297303),
    @(-0.25945945945945947, 0.16216216216216217),
    @(-0.22702702702702704, 0.16216216216216217),
    @(-0.22702702702702704, 0.097297297297297303),
    @(-0.22702702702702704, 0.064864864864864868),
    @(-0.22702702702702704, 0),
    @(-0.25945945945945947, 0),
    @(-0.25945945945945947, 0.064864864864864868),
    @(-0.25945945945945947, 0.097297297297297303)
)
for ($i = 0; $i -lt $efVals2.Length; $i++) {
    $row = 79 + $i
    $pair = $efVals2[$i]
    $ws.Cells.Item($row, 5).Value2 = $pair[0]
    $ws.Cells.Item($row, 6).Value2 = $pair[1]
}

# --- 5) Formulas: standalone (row 69, row 79) ---
$ws.Range("G69").Formula = "=10000*E69"
$ws.Range("H69").Formula = "=10000*F69"
$ws.Range("G79").Formula = "=10000*E79"
$ws.Range("H79").Formula = "=10000*F79"

# --- 6) Formulas: shared ranges G70:G77 / H70:H77 ---
$ws.Range("G70:G77").Formula = "=10000*E70"
$ws.Range("H70:H77").Formula = "=10000*F70"

# --- 7) Formulas: shared ranges G80:G87 / H80:H87 ---
$ws.Range("G80:G87").Formula = "=10000*E80"
$ws.Range("H80:H87").Formula = "=10000*F80"

# --- 8) Update selection/view to match the edited region ---
$ws.Range("J85").Select()
